# Apply the "Output auxiliary files in fmudesign init" edit:
#  - remove the explanatory cell comments (author decided the example
#    workbook should ship lean, without the design-notes balloons)
#  - drop the extra illustrative default-value rows (FAULT_POSITION,
#    OWC1/2/3, MULTZ_ILE, PARAM1-4) from the "defaultvalues" sheet,
#    keeping only RMS_SEED / VEL_MODEL / COHIBA_MODE
#  - restore the B4 centered style on "general_input" that Excel
#    re-applies once the sheet is touched
#  - leave designinput as the active/selected sheet, matching the
#    resaved workbook view

$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("general_input")
$wsDesign  = $wb.Worksheets.Item("designinput")
$wsDefault = $wb.Worksheets.Item("defaultvalues")

# --- remove all cell comments on every sheet ---
$wsGeneral.Range("B1").Comment.Delete()
$wsGeneral.Range("B2").Comment.Delete()
$wsGeneral.Range("B3").Comment.Delete()
$wsGeneral.Range("B4").Comment.Delete()

$wsDesign.Range("B1").Comment.Delete()
$wsDesign.Range("E1").Comment.Delete()
$wsDesign.Range("I1").Comment.Delete()

$wsDefault.Range("B1").Comment.Delete()

# --- trim "defaultvalues" down to param_name/default_value, RMS_SEED,
#     VEL_MODEL and COHIBA_MODE only ---
$wsDefault.Range("A5:B12").EntireRow.Delete()
$wsDefault.Range("A3:B3").EntireRow.Delete()

# --- re-apply the centered alignment style on general_input!B4 ---
$wsGeneral.Range("B4").HorizontalAlignment = -4108

# --- selection / active sheet bookkeeping to match the resaved file ---
$wsGeneral.Range("C3").Select()
$wsDesign.Range("E18").Select()
$wsDefault.Range("C12").Select()
$wsDesign.Activate()
